$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from 45203 (2023-10-04)
# to 45205 (2023-10-06) for every data row (rows 2 through 499).
$ws.Range("C2:C499").Value = 45205
